$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.05
$ws.Range("L2").Value = 5
$ws.Range("W2").Value = 9.5
$ws.Range("Z2").Value = 12
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 34
$ws.Range("AJ2").Value = 19
$ws.Range("AL2").Value = 41
$ws.Range("AM2").Value = 41
$ws.Range("AO2").Value = 7.5
$ws.Range("AQ2").Value = 21
$ws.Range("AR2").Value = 34
$ws.Range("AU2").Value = 7.5
$ws.Range("AW2").Value = 351
$ws.Range("AX2").Value = 7.5
$ws.Range("AY2").Value = 26
$ws.Range("AZ2").Value = 26
$ws.Range("BA2").Value = 81
$ws.Range("BB2").Value = 81
$ws.Range("BC2").Value = 151
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.3
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 1.95
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 2.1
$ws.Range("Y3").Value = 9
$ws.Range("AG3").Value = 151
$ws.Range("AL3").Value = 26
$ws.Range("G4").Value = 2.05
$ws.Range("I4").Value = 3.75
$ws.Range("L4").Value = 4.33
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("Z4").Value = 19
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 8
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 17
$ws.Range("AO4").Value = 12
$ws.Range("AT4").Value = 2.5
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 3.8
$ws.Range("L5").Value = 4.5
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("X5").Value = 9
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 23
$ws.Range("AC5").Value = 5.5
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 17
$ws.Range("AQ5").Value = 51
$ws.Range("G6").Value = 2.05
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 2.88
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.38
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 8.5
$ws.Range("AH6").Value = 9
$ws.Range("AT6").Value = 2.38
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 67
$ws.Range("AZ6").Value = 34
$ws.Range("G7").Value = 2.15
$ws.Range("I7").Value = 3.5
$ws.Range("K7").Value = 2.2
$ws.Range("L7").Value = 3.75
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 3.75
$ws.Range("Q7").Value = 1.83
$ws.Range("R7").Value = 2.03
$ws.Range("S7").Value = 1.36
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.1
$ws.Range("AC7").Value = 11
$ws.Range("AI7").Value = 19
$ws.Range("AL7").Value = 26
$ws.Range("AP7").Value = 21
$ws.Range("AT7").Value = 3
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 1.44
$ws.Range("AD8").Value = 9
$ws.Range("AE8").Value = 29
$ws.Range("AH8").Value = 5
$ws.Range("AI8").Value = 5.5
$ws.Range("AJ8").Value = 9.5
$ws.Range("AU8").Value = 11
$ws.Range("BC8").Value = 251
